$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 13 - this shifts the existing row 13
# (and everything below it) down by one row, so the former row 13
# becomes row 14, former row 106 becomes row 107, etc.
$ws.Rows.Item(13).Insert()

# Populate the newly inserted row 13 with the new weekly record.
$ws.Cells.Item(13, 1).Value = 11
$ws.Cells.Item(13, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(13, 3).Value = "Bíobío"
$ws.Cells.Item(13, 4).Value = 44635
$ws.Cells.Item(13, 5).Value = 8
$ws.Cells.Item(13, 6).Value = 100112032
$ws.Cells.Item(13, 7).Value = "Zapallo italiano"
$ws.Cells.Item(13, 8).Value = "Sin especificar"
$ws.Cells.Item(13, 9).Value = "Primera"
$ws.Cells.Item(13, 10).Value = 220
$ws.Cells.Item(13, 11).Value = 12000
$ws.Cells.Item(13, 12).Value = 13000
$ws.Cells.Item(13, 13).Value = 12545
$ws.Cells.Item(13, 14).Value = "$/caja 60 unidades"
$ws.Cells.Item(13, 15).Value = "Región Metropolitana"
$ws.Cells.Item(13, 16).Value = 209
$ws.Cells.Item(13, 17).Value = 60
$ws.Cells.Item(13, 18).Value = "Hortaliza"
